$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.201.89'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '2.321.84'
$ws.Range("E3").Value = '  +1.02%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.69%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +2.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.03'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.42%  '
$ws.Range("E12").Value = '  -0.70%  '
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("E14").Value = '  +2.29%  '
$ws.Range("D15").Value = '2.683.32'
$ws.Range("E15").Value = '  +1.06%  '
$ws.Range("D16").Value = '2.337.45'
$ws.Range("E16").Value = '  +1.08%  '
$ws.Range("E17").Value = '  -1.25%  '
$ws.Range("D18").Value = '43.099.72'
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.97%  '
$ws.Range("E21").Value = '  +1.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.19'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.13%  '
$ws.Range("E24").Value = '  -1.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  +3.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.56%  '
$ws.Range("E31").Value = '  -6.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("E34").Value = '  +5.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.66'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.83%  '
$ws.Range("E36").Value = '  -0.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0699'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.11%  '
$ws.Range("E38").Value = '  +2.59%  '
$ws.Range("E39").Value = '  +0.49%  '
$ws.Range("E40").Value = '  -0.50%  '
$ws.Range("E41").Value = '  +0.50%  '
$ws.Range("D42").Value = '1.996.34'
$ws.Range("E43").Value = '  +1.62%  '
$ws.Range("E44").Value = '  -4.44%  '
$ws.Range("E45").Value = '  +1.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.86'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '76.69'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.11'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.13%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.548.46'
$ws.Range("E50").Value = '  +0.72%  '
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.86'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +11.97%  '
